# Updated cryptos list on Sat Mar 25 17:33:24 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.569.64"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "1.752.15"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'324.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").Value = "'0.4570"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.19%  "
$ws.Range("D8").Value = "'0.3580"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.90%  "
$ws.Range("D9").Value = "'0.07506"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.88%  "
$ws.Range("D10").Value = "'42.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.58%  "
$ws.Range("D11").Value = "'1.095"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "'20.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").Value = "'6.003"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'7.085"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.14%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.754.86"
$ws.Range("E16").Value = "  -1.76%  "
$ws.Range("D17").Value = "'92.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").Value = "'0.00001063"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").Value = "'0.06418"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").Value = "'1.003"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'16.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("D22").Value = "'5.818"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.54%  "
$ws.Range("D23").Value = "27.647.71"
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("D24").Value = "'11.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("D25").Value = "'2.109"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("D26").Value = "'164.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.21%  "
$ws.Range("D27").Value = "'20.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("D28").Value = "1.957.24"
$ws.Range("E28").Value = "  -1.34%  "
$ws.Range("D29").Value = "'2.076"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.85%  "
$ws.Range("D30").Value = "'126.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("E31").Value = "  -6.74%  "
$ws.Range("D32").Value = "'0.09175"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.36%  "
$ws.Range("D33").Value = "'3.669"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.35%  "
$ws.Range("D34").Value = "'5.532"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.28%  "
$ws.Range("D35").Value = "'11.86"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.08%  "
$ws.Range("D36").Value = "'0.02296"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").Value = "'0.2100"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("D38").Value = "'0.06047"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("D39").Value = "'4.973"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("D40").Value = "'0.6332"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("D41").Value = "'1.209"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.44%  "
$ws.Range("D42").Value = "'1.380"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.71%  "
$ws.Range("D43").Value = "'7.787"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("D44").Value = "'13.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.86%  "
$ws.Range("D45").Value = "'0.5905"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("D46").Value = "'3.716"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("D47").Value = "'123.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.63%  "
$ws.Range("D48").Value = "'1.944"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.48%  "
$ws.Range("D49").Value = "'1.144"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.95%  "
$ws.Range("D50").Value = "'0.06865"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.59%  "
$ws.Range("D51").Value = "'72.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.72%  "
